$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 15: Meagan / mmcleod@illumina.com / Green (favorite color email automation)
$ws.Range("A15").Value = "Meagan"
$ws.Range("B15").Value = "mmcleod@illumina.com"
$ws.Range("C15").Value = "Green"

# Style the new row with a bordered, wrap-text look (new cellXf: font+border+wrapText)
$newRowRange = $ws.Range("A15:C15")
$newRowRange.Borders.Weight = -4138
$newRowRange.Borders.Color = 13421772
$newRowRange.WrapText = $true

# Row 16: a lone helper value in I16
$ws.Cells.Item(16, 9).Value = 4
